$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-09-13"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = "AA"
$ws.Range("C8").Value = "44CDX12"
$ws.Range("D8").Value = "MAM "

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-08-13"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = "JJ"
$ws.Range("C9").Value = "456CDX176"
$ws.Range("D9").Value = "AMM "
